# Applies the numeric-value updates from the scheduled-runner sheet refresh.
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 979.8
$ws.Range("I2").Value = 1400
$ws.Range("K2").Value = 1400
$ws.Range("M2").Value = -1287
$ws.Range("H5").Value = 1704
$ws.Range("I5").Value = 70
$ws.Range("J5").Value = 2521
$ws.Range("K5").Value = 70
$ws.Range("L5").Value = 2521
$ws.Range("M5").Value = 45
$ws.Range("N5").Value = -2751
$ws.Range("H121").Value = 4026.7273
$ws.Range("J121").Value = 4026.7273
$ws.Range("L121").Value = 12080.1819
$ws.Range("N121").Value = -15574.1819

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 21278996
$ws.Range("I74").Value = 58825180
$ws.Range("K74").Value = 58825180
$ws.Range("M74").Value = -58824306
$ws.Range("H77").Value = 21278996
$ws.Range("I77").Value = 58825180
$ws.Range("K77").Value = 294125900
$ws.Range("M77").Value = -294121532
$ws.Range("H132").Value = 20863834
$ws.Range("I132").Value = 1433.7
$ws.Range("K132").Value = 4301.1
$ws.Range("M132").Value = -1771.1

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H74").Value = 37922.2
$ws.Range("I74").Value = 49994
$ws.Range("J74").Value = 34904.25
$ws.Range("K74").Value = 49994
$ws.Range("L74").Value = 34904.25
$ws.Range("M74").Value = -49058
$ws.Range("N74").Value = -36776.25
$ws.Range("H77").Value = 37922.2
$ws.Range("I77").Value = 49994
$ws.Range("J77").Value = 34904.25
$ws.Range("K77").Value = 149982
$ws.Range("L77").Value = 104712.75
$ws.Range("M77").Value = -145302
$ws.Range("N77").Value = -114072.75
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 6350.2354
$ws.Range("I22").Value = 9398.091
$ws.Range("J22").Value = 762.5
$ws.Range("K22").Value = 9398.091
$ws.Range("L22").Value = 762.5
$ws.Range("M22").Value = -9048.091
$ws.Range("N22").Value = -1462.5
$ws.Range("H31").Value = 73535580
$ws.Range("I31").Value = 4168.5
$ws.Range("J31").Value = 178580460
$ws.Range("K31").Value = 4168.5
$ws.Range("L31").Value = 178580460
$ws.Range("M31").Value = -3873.5
$ws.Range("N31").Value = -178581050
$ws.Range("H34").Value = 73535580
$ws.Range("I34").Value = 4168.5
$ws.Range("J34").Value = 178580460
$ws.Range("K34").Value = 4168.5
$ws.Range("L34").Value = 178580460
$ws.Range("M34").Value = -3966.5
$ws.Range("N34").Value = -178580864
$ws.Range("H48").Value = 37708.25
$ws.Range("J48").Value = 37708.25
$ws.Range("L48").Value = 37708.25
$ws.Range("N48").Value = -38660.25
$ws.Range("H57").Value = 39098.4
$ws.Range("J57").Value = 41373
$ws.Range("L57").Value = 41373
$ws.Range("N57").Value = -42493
$ws.Range("H99").Value = 7425.5713
$ws.Range("I99").Value = 7080
$ws.Range("J99").Value = 9499
$ws.Range("K99").Value = 7080
$ws.Range("L99").Value = 9499
$ws.Range("M99").Value = -5582
$ws.Range("N99").Value = -12495
$ws.Range("H122").Value = 4051220.2
$ws.Range("I122").Value = 1942.375
$ws.Range("K122").Value = 5827.125
$ws.Range("M122").Value = -3377.125
$ws.Range("H126").Value = 7425.5713
$ws.Range("I126").Value = 7080
$ws.Range("J126").Value = 9499
$ws.Range("K126").Value = 21240
$ws.Range("L126").Value = 28497
$ws.Range("M126").Value = -18770
$ws.Range("N126").Value = -33437

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 857.0952
$ws.Range("J34").Value = 1436.2727
$ws.Range("L34").Value = 4308.8181
$ws.Range("N34").Value = -4476.8181
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H55").Value = 6064492
$ws.Range("I55").Value = 2719.6
$ws.Range("J55").Value = 11115969
$ws.Range("K55").Value = 8158.799999999999
$ws.Range("L55").Value = 33347907
$ws.Range("M55").Value = -7981.799999999999
$ws.Range("N55").Value = -33348261
$ws.Range("H131").Value = 31819424
$ws.Range("J131").Value = 1719.5
$ws.Range("L131").Value = 5158.5
$ws.Range("N131").Value = -15238.5
$ws.Range("H132").Value = 5560235.5
$ws.Range("J132").Value = 9530661
$ws.Range("L132").Value = 85775949
$ws.Range("N132").Value = -85781009
$ws.Range("H134").Value = 2928.9
$ws.Range("I134").Value = 1032.1111
$ws.Range("J134").Value = 20000
$ws.Range("K134").Value = 3096.3333
$ws.Range("L134").Value = 60000
$ws.Range("M134").Value = 1973.6667
$ws.Range("N134").Value = -70140
$ws.Range("H139").Value = 2082.8333
$ws.Range("I139").Value = 2149.4707
$ws.Range("J139").Value = 950
$ws.Range("K139").Value = 6448.4121
$ws.Range("L139").Value = 2850
$ws.Range("M139").Value = -1308.4121
$ws.Range("N139").Value = -13130
$ws.Range("H141").Value = 14149.286
$ws.Range("I141").Value = 9168.833000000001
$ws.Range("J141").Value = 44032
$ws.Range("K141").Value = 27506.499
$ws.Range("L141").Value = 132096
$ws.Range("M141").Value = -22326.499
$ws.Range("N141").Value = -142456

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3674.9697
$ws.Range("I113").Value = 3235.1667
$ws.Range("K113").Value = 3235.1667
$ws.Range("M113").Value = -1065.1667
$ws.Range("H122").Value = 38465030
$ws.Range("I122").Value = 2676.125
$ws.Range("J122").Value = 100004800
$ws.Range("K122").Value = 8028.375
$ws.Range("L122").Value = 300014400
$ws.Range("M122").Value = -5578.375
$ws.Range("N122").Value = -300019300
$ws.Range("H126").Value = 21505536
$ws.Range("I126").Value = 14433928
$ws.Range("K126").Value = 43301784
$ws.Range("M126").Value = -43299314

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 333352160
$ws.Range("I132").Value = 4000
$ws.Range("K132").Value = 12000
$ws.Range("M132").Value = -9470

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 8798.933999999999
$ws.Range("I126").Value = 8641.786
$ws.Range("J126").Value = 10999
$ws.Range("K126").Value = 25925.358
$ws.Range("L126").Value = 32997
$ws.Range("M126").Value = -23455.358
$ws.Range("N126").Value = -37937
$ws.Range("H132").Value = 2028.8572
$ws.Range("I132").Value = 2850
$ws.Range("J132").Value = 1700.4
$ws.Range("K132").Value = 8550
$ws.Range("L132").Value = 5101.200000000001
$ws.Range("M132").Value = -6020
$ws.Range("N132").Value = -10161.2
$ws.Range("H136").Value = 1617.8
$ws.Range("I136").Value = 1617.8
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4853.4
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -2303.4
$ws.Range("N136").ClearContents()
$ws.Range("H140").Value = 244000
$ws.Range("J140").Value = 244000
$ws.Range("L140").Value = 244000
$ws.Range("N140").Value = -254360
